# Updated cryptos list values (price / 1h volume change) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '38.702.01'
$ws.Range('E2').Value = '  +2.47%  '
$ws.Range('D3').Value = '2.088.33'
$ws.Range('E3').Value = '  +2.29%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = "'228.22"
$ws.Range('E5').Value = '  +0.31%  '
$ws.Range('D6').Value = "'0.614"
$ws.Range('E6').Value = '  +0.99%  '
$ws.Range('D7').Value = "'60.65"
$ws.Range('E7').Value = '  +1.55%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').Value = "'0.384"
$ws.Range('E9').Value = '  +2.20%  '
$ws.Range('D10').Value = "'0.0836"
$ws.Range('E10').Value = '  -0.88%  '
$ws.Range('D11').Value = "'0.103"
$ws.Range('E11').Value = '  -0.51%  '
$ws.Range('D12').Value = '2.398.34'
$ws.Range('E12').Value = '  +2.43%  '
$ws.Range('D13').Value = "'14.94"
$ws.Range('E13').Value = '  +3.66%  '
$ws.Range('E14').Value = '  +4.01%  '
$ws.Range('D15').Value = "'0.795"
$ws.Range('E15').Value = '  +3.33%  '
$ws.Range('D16').Value = "'5.47"
$ws.Range('E16').Value = '  +0.11%  '
$ws.Range('D17').Value = '2.089.66'
$ws.Range('E17').Value = '  +1.96%  '
$ws.Range('D18').Value = '38.705.88'
$ws.Range('E18').Value = '  +2.65%  '
$ws.Range('D19').Value = "'71.65"
$ws.Range('E19').Value = '  +3.20%  '
$ws.Range('D20').Value = "'6.02"
$ws.Range('E20').Value = '  +1.83%  '
$ws.Range('E21').Value = '  +1.32%  '
$ws.Range('D22').Value = "'226.45"
$ws.Range('E22').Value = '  +1.31%  '
$ws.Range('D24').Value = "'2.44"
$ws.Range('E24').Value = '  +0.24%  '
$ws.Range('E25').Value = '  +2.88%  '
$ws.Range('D26').Value = "'170.73"
$ws.Range('E26').Value = '  +0.89%  '
$ws.Range('D27').Value = "'9.43"
$ws.Range('E27').Value = '  +0.79%  '
$ws.Range('E28').Value = '  +7.03%  '
$ws.Range('D29').Value = "'1.44"
$ws.Range('E29').Value = '  +11.70%  '
$ws.Range('D30').Value = "'19.16"
$ws.Range('E30').Value = '  +2.01%  '
$ws.Range('E31').Value = '  +0.95%  '
$ws.Range('E32').Value = '  +4.21%  '
$ws.Range('D33').Value = "'4.50"
$ws.Range('E33').Value = '  +2.77%  '
$ws.Range('D34').Value = "'4.70"
$ws.Range('E34').Value = '  +4.67%  '
$ws.Range('E35').Value = '  +2.09%  '
$ws.Range('E36').Value = '  +1.79%  '
$ws.Range('D37').Value = "'6.39"
$ws.Range('E37').Value = '  -2.26%  '
$ws.Range('D38').Value = "'3.53"
$ws.Range('E38').Value = '  +2.47%  '
$ws.Range('E39').Value = '  +0.14%  '
$ws.Range('D40').Value = "'18.25"
$ws.Range('E40').Value = '  +1.45%  '
$ws.Range('D41').Value = '1.538.79'
$ws.Range('E41').Value = '  +0.66%  '
$ws.Range('E42').Value = '  +3.22%  '
$ws.Range('E43').Value = '  +3.53%  '
$ws.Range('B44').Value = 'HuobiToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D44').Value = "'2.82"
$ws.Range('E44').Value = '  -0.76%  '
$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D45').Value = "'0.0924"
$ws.Range('E45').Value = '  +2.12%  '
$ws.Range('D46').Value = "'7.69"
$ws.Range('E46').Value = '  +8.39%  '
$ws.Range('E47').Value = '  +0.75%  '
$ws.Range('D48').Value = "'4.11"
$ws.Range('E48').Value = '  -1.91%  '
$ws.Range('E49').Value = '  +2.26%  '
$ws.Range('E50').Value = '  +1.04%  '
$ws.Range('D51').Value = '2.289.08'
$ws.Range('E51').Value = '  +2.60%  '

Write-Host "Updated crypto values"
